$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Configuration")

# --- Update the "Configuration" sheet content ---
# Row 6: rename "Filename" -> "FilenamePattern" and add a description in column C.
$wsConfig.Range("A6").Value = "FilenamePattern"
$wsConfig.Range("C6").Value = 'Here you may define the file name pattern of the generated files. You can use all Variables of the template, e. g. ''letter-''${receiver}''.'

# --- Column widths on the "Configuration" sheet ---
$wsConfig.Columns.Item(1).ColumnWidth = 12.830729166666666
$wsConfig.Columns.Item(3).ColumnWidth = 99.16666666666667

# --- Switch the active tab from "Variables" to "Configuration" and update the selection ---
$wsConfig.Activate()
$wsConfig.Range("A1:C1").Select()
